# Revert "PowerPoint writer: consolidate text run nodes."
#
# The caption textbox ("TextBox 3" on slide 1) currently stores its text
# as 4 runs: "Followed ", "by ", "a ", "picture". The target OOXML keeps
# the words and the separating spaces as distinct runs instead: "Followed",
# " ", "by", " ", "a", " ", "picture" (7 runs total), with the overall
# text content unchanged.
#
# Assigning `.Text` on a `Characters(start, length)` sub-range that covers
# exactly the prefix of an existing run (and re-supplying that same
# substring) splits that run in two at the boundary without touching the
# remainder of the text - which is exactly the run-boundary change the
# diff describes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("TextBox 3")
$tr = $sh.TextFrame.TextRange

# Sanity check: unmodified caption text before the split.
# "Followed " + "by " + "a " + "picture" = "Followed by a picture"

# Split run 1 "Followed " (chars 1-9) into "Followed" (1-8) + " " (9)
$tr.Characters(1, 8).Text = "Followed"

# Split run 2 "by " (chars 10-12) into "by" (10-11) + " " (12)
$tr.Characters(10, 2).Text = "by"

# Split run 3 "a " (chars 13-14) into "a" (13) + " " (14)
$tr.Characters(13, 1).Text = "a"

# Run 4 "picture" (chars 15-21) is left untouched.
